# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# Values that read as pure numbers (e.g. "1.00", "580.32") are written with a
# leading apostrophe so Excel keeps them as literal text, matching the
# original sheet where every data cell is stored as text (inlineStr).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.279.75"
$ws.Range("E2").Value = "  -5.44%  "
$ws.Range("D3").Value = "3.208.26"
$ws.Range("E3").Value = "  -8.41%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'580.32"
$ws.Range("E5").Value = "  -3.69%  "
$ws.Range("D6").Value = "'151.27"
$ws.Range("E6").Value = "  -13.79%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "3.204.40"
$ws.Range("E8").Value = "  -8.44%  "
$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  -13.37%  "
$ws.Range("D10").Value = "'0.170"
$ws.Range("E10").Value = "  -11.49%  "
$ws.Range("D11").Value = "'5.99"
$ws.Range("E11").Value = "  -17.24%  "
$ws.Range("D12").Value = "'0.469"
$ws.Range("E12").Value = "  -19.15%  "
$ws.Range("D13").Value = "'37.60"
$ws.Range("E13").Value = "  -18.55%  "
$ws.Range("D14").Value = "'0.0000232"
$ws.Range("E14").Value = "  -15.22%  "
$ws.Range("D15").Value = "3.732.67"
$ws.Range("E15").Value = "  -8.45%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "66.228.90"
$ws.Range("E16").Value = "  -5.70%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.218.16"
$ws.Range("E17").Value = "  -8.16%  "
$ws.Range("E18").Value = "  -4.04%  "
$ws.Range("D19").Value = "'518.26"
$ws.Range("E19").Value = "  -15.11%  "
$ws.Range("D20").Value = "'6.78"
$ws.Range("E20").Value = "  -18.04%  "
$ws.Range("D21").Value = "'14.49"
$ws.Range("E21").Value = "  -16.34%  "
$ws.Range("D22").Value = "'0.740"
$ws.Range("E22").Value = "  -15.39%  "
$ws.Range("D23").Value = "'7.44"
$ws.Range("E23").Value = "  -16.92%  "
$ws.Range("D24").Value = "'83.27"
$ws.Range("E24").Value = "  -14.69%  "
$ws.Range("D25").Value = "'13.05"
$ws.Range("E25").Value = "  -15.79%  "
$ws.Range("D26").Value = "'0.988"
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("D27").Value = "'3.24"
$ws.Range("E27").Value = "  -12.71%  "
$ws.Range("D28").Value = "'28.16"
$ws.Range("E28").Value = "  -16.39%  "
$ws.Range("D29").Value = "'2.05"
$ws.Range("E29").Value = "  -19.22%  "
$ws.Range("D30").Value = "'7.25"
$ws.Range("E30").Value = "  -19.09%  "
$ws.Range("E31").Value = "  -12.70%  "
$ws.Range("D32").Value = "'2.45"
$ws.Range("E32").Value = "  -17.05%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.38%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").Value = "'516.78"
$ws.Range("E34").Value = "  -18.57%  "
$ws.Range("D35").Value = "'5.57"
$ws.Range("E35").Value = "  -18.32%  "
$ws.Range("D36").Value = "'6.21"
$ws.Range("E36").Value = "  -22.19%  "
$ws.Range("D37").Value = "'53.43"
$ws.Range("E37").Value = "  -5.69%  "
$ws.Range("D38").Value = "'0.0415"
$ws.Range("E38").Value = "  -11.51%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.0826"
$ws.Range("E39").Value = "  -16.44%  "
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").Value = "'8.95"
$ws.Range("E40").Value = "  -16.43%  "
$ws.Range("D41").Value = "'0.117"
$ws.Range("E41").Value = "  -17.20%  "
$ws.Range("D42").Value = "2.813.77"
$ws.Range("E42").Value = "  -16.08%  "
$ws.Range("D43").Value = "'2.60"
$ws.Range("E43").Value = "  -26.53%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.249"
$ws.Range("E45").Value = "  -18.83%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").Value = "'2.35"
$ws.Range("E46").Value = "  -17.97%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'25.22"
$ws.Range("E47").Value = "  -21.41%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "'122.67"
$ws.Range("E48").Value = "  -8.45%  "
$ws.Range("B49").Value = "PEPE"
$ws.Range("C49").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D49").Value = "0.0₃0531"
$ws.Range("E49").Value = "  -27.58%  "
$ws.Range("E50").Value = "  -14.55%  "
$ws.Range("E51").Value = "  -21.26%  "
